# Update the three "Email Id" values in the sample import sheet so the
# fixture can be reused for the export-template test cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "test0100@gmail.com"
$ws.Range("E3").Value = "test0200@gmail.com"
$ws.Range("E4").Value = "test0300@gmail.com"

# Leave the selection on the last edited cell, matching the author's saved view.
$ws.Range("E4").Select()
